$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (counts)
$ws.Range("B2").Value = 1659
$ws.Range("B3").Value = 1187
$ws.Range("B4").Value = 971
$ws.Range("B5").Value = 431

# Update column A values (group ids) - reorder rows 3,4,5
$ws.Range("A3").Value = 22
$ws.Range("A4").Value = 21
$ws.Range("A5").Value = 12
